# Insert a new price record for Berenjena (Macroferia Regional de Talca)
# right above the existing row 52, shifting all subsequent rows down by one.
# The sheet's dimension grows from A1:R92 to A1:R93 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 44589
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 100112001
$ws.Range("G52").Value = "Berenjena"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 200
$ws.Range("K52").Value = 7000
$ws.Range("L52").Value = 7000
$ws.Range("M52").Value = 7000
$ws.Range("N52").Value = "$/caja 50 unidades"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 140
$ws.Range("Q52").Value = 50
$ws.Range("R52").Value = "Hortaliza"
